# Improve testcases + format
# Adds "Events" / "Mean VisibilityKm" columns (K/L) with some sample data,
# clears a stray value in C4, widens the narrow columns to match the
# wide "Mean VisibilityKm" header column, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells -------------------------------------------------
$ws.Range("K1").Value = "Events"
$ws.Range("L1").Value = "Mean VisibilityKm"

# --- New data cells -----------------------------------------------------
$ws.Range("K2").Value = "Rain"

$ws.Range("L3").Value = 12.8

$ws.Range("K4").Value = "Thunder"
$ws.Range("L4").Value = 11.98

# --- Remove the stray duplicate value in C4 ------------------------------
$ws.Range("C4").ClearContents()

# --- Column widths: narrow A:G and I:L columns to match the rest --------
$ws.Range("A1:G4").ColumnWidth = 7.836666666666667
$ws.Range("I1:L4").ColumnWidth = 7.836666666666667

# --- Move the active selection to C4 -------------------------------------
[void]$ws.Range("C4").Select()
